$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Firstname value on row 2 from "ganesh" to "ravi"
$ws.Range("C2").Value = "ravi"

# Delete the entire row 3 (the "ravisankar" test case), shifting cells up
$ws.Rows(3).Delete()

# Update selection to reflect the after-state (C3, where row 3 used to be)
$ws.Range("C3").Select()
